$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.086.69'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '1.651.91'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').Value = '''218.48'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').Value = '''0.5251'
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D8').Value = '''0.2673'
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('D9').Value = '''0.06365'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '''20.52'
$ws.Range('E10').Value = '  -1.97%  '
$ws.Range('D11').Value = '''0.07684'
$ws.Range('E11').Value = '  -1.87%  '
$ws.Range('D12').Value = '''4.592'
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.654.82'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.880.10'
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '''0.5613'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '0.0₅8226'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('D17').Value = '''65.40'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '26.095.46'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '''4.683'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').Value = '''10.34'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').Value = '''190.90'
$ws.Range('E22').Value = '  -4.09%  '
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('D25').Value = '''145.90'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('D26').Value = '''0.1200'
$ws.Range('E26').Value = '  -1.02%  '
$ws.Range('D27').Value = '''7.239'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '''15.95'
$ws.Range('E28').Value = '  -1.56%  '
$ws.Range('D29').Value = '''1.496'
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('D30').Value = '''0.05643'
$ws.Range('E30').Value = '  -4.01%  '
$ws.Range('D31').Value = '''1.269'
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').Value = '''3.491'
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('D33').Value = '''3.375'
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('D34').Value = '''1.577'
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('D35').Value = '''2.794'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '''2.409'
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '''0.9444'
$ws.Range('E37').Value = '  -1.81%  '
$ws.Range('D38').Value = '''0.5769'
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('D39').Value = '''0.01590'
$ws.Range('E39').Value = '  -1.67%  '
$ws.Range('D40').Value = '''5.968'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('E41').Value = '  -0.53%  '
$ws.Range('D42').Value = '''0.8407'
$ws.Range('E42').Value = '  -1.83%  '
$ws.Range('D43').Value = '1.021.14'
$ws.Range('E43').Value = '  -5.15%  '
$ws.Range('D44').Value = '''101.45'
$ws.Range('E44').Value = '  -1.28%  '
$ws.Range('D45').Value = '1.791.28'
$ws.Range('E45').Value = '  -0.82%  '
$ws.Range('D46').Value = '''58.47'
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₈105'
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.05350'
$ws.Range('E48').Value = '  +3.95%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').Value = '''1.003'
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''8.049'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '''0.4343'
$ws.Range('E51').Value = '  -1.63%  '
